# Auto-generated edit script: updates crypto price/volume figures
# and fixes the swapped TrustWalletToken / InjectiveProtocol rows (43-44).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing it to stay a text cell (matches the
# original inline-string cells) instead of letting Excel auto-coerce
# numeric-looking strings (e.g. "17.70", "1.22") into real numbers, and
# without leaving a stray text-format style behind on the cell.
function Set-TextCell($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" '37.421.37'
Set-TextCell $ws "E2" '  +5.06%  '
Set-TextCell $ws "D3" '2.053.08'
Set-TextCell $ws "E3" '  +3.50%  '
Set-TextCell $ws "E4" '  +0.02%  '
Set-TextCell $ws "D5" '252.75'
Set-TextCell $ws "E5" '  +3.09%  '
Set-TextCell $ws "D6" '0.653'
Set-TextCell $ws "E6" '  +2.48%  '
Set-TextCell $ws "D7" '65.46'
Set-TextCell $ws "E7" '  +13.41%  '
Set-TextCell $ws "E8" '  -0.01%  '
Set-TextCell $ws "E9" '  +6.43%  '
Set-TextCell $ws "D10" '59.25'
Set-TextCell $ws "E10" '  +1.36%  '
Set-TextCell $ws "E11" '  +4.61%  '
Set-TextCell $ws "E12" '  +1.46%  '
Set-TextCell $ws "D13" '0.928'
Set-TextCell $ws "E13" '  -1.64%  '
Set-TextCell $ws "D14" '14.84'
Set-TextCell $ws "E14" '  +2.72%  '
Set-TextCell $ws "E15" '  +25.66%  '
Set-TextCell $ws "D16" '2.354.00'
Set-TextCell $ws "E16" '  +3.56%  '
Set-TextCell $ws "E17" '  +5.52%  '
Set-TextCell $ws "D18" '2.057.08'
Set-TextCell $ws "E18" '  +3.81%  '
Set-TextCell $ws "D19" '37.316.46'
Set-TextCell $ws "E19" '  +4.82%  '
Set-TextCell $ws "D20" '73.51'
Set-TextCell $ws "E20" '  +3.11%  '
Set-TextCell $ws "D21" '0.0₃0875'
Set-TextCell $ws "E21" '  +3.84%  '
Set-TextCell $ws "D22" '5.49'
Set-TextCell $ws "E22" '  +6.30%  '
Set-TextCell $ws "D23" '239.86'
Set-TextCell $ws "E23" '  +3.04%  '
Set-TextCell $ws "D24" '2.68'
Set-TextCell $ws "E24" '  +4.95%  '
Set-TextCell $ws "E25" '  +0.07%  '
Set-TextCell $ws "E26" '  +5.10%  '
Set-TextCell $ws "D27" '10.05'
Set-TextCell $ws "E27" '  +9.80%  '
Set-TextCell $ws "D28" '161.92'
Set-TextCell $ws "E28" '  -1.68%  '
Set-TextCell $ws "D29" '19.99'
Set-TextCell $ws "E29" '  +4.14%  '
Set-TextCell $ws "D30" '0.123'
Set-TextCell $ws "E30" '  +29.16%  '
Set-TextCell $ws "E31" '  +8.42%  '
Set-TextCell $ws "E32" '  +2.74%  '
Set-TextCell $ws "D33" '1.22'
Set-TextCell $ws "E33" '  +9.68%  '
Set-TextCell $ws "D34" '4.69'
Set-TextCell $ws "E34" '  +8.43%  '
Set-TextCell $ws "D35" '0.0627'
Set-TextCell $ws "E35" '  +5.85%  '
Set-TextCell $ws "D36" '2.45'
Set-TextCell $ws "E36" '  +1.63%  '
Set-TextCell $ws "E37" '  +0.00%  '
Set-TextCell $ws "E38" '  +3.99%  '
Set-TextCell $ws "E39" '  +15.04%  '
Set-TextCell $ws "D40" '3.05'
Set-TextCell $ws "E40" '  +36.33%  '
Set-TextCell $ws "E41" '  +17.16%  '
Set-TextCell $ws "E42" '  +4.65%  '
Set-TextCell $ws "B43" 'TrustWalletToken'
Set-TextCell $ws "C43" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws "D43" '1.24'
Set-TextCell $ws "E43" '  +1.99%  '
Set-TextCell $ws "B44" 'InjectiveProtocol'
Set-TextCell $ws "C44" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws "D44" '17.70'
Set-TextCell $ws "E44" '  +10.15%  '
Set-TextCell $ws "E45" '  +5.94%  '
Set-TextCell $ws "E46" '  +3.07%  '
Set-TextCell $ws "D47" '96.67'
Set-TextCell $ws "E47" '  +5.15%  '
Set-TextCell $ws "D48" '7.91'
Set-TextCell $ws "E48" '  +2.98%  '
Set-TextCell $ws "D49" '1.417.06'
Set-TextCell $ws "E49" '  +3.45%  '
Set-TextCell $ws "E50" '  +1.95%  '
Set-TextCell $ws "D51" '46.62'
Set-TextCell $ws "E51" '  -0.07%  '
